# SPT-327: replace the hard-coded signature date ("1st October 2024")
# in the SigBlock/SigDate paragraph with a plain "Date" placeholder run.
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("1st October 2024", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "Date", 2)

if (-not $found) {
    throw "edit.ps1: could not find the '1st October 2024' signature date text to replace"
}
